# Weekly fruit/vegetable price update: rotate the D/L/M/N/O/P/S values
# across rows 8, 11, 15, 9, 12, 16, 10 (cyclic shift, each row takes on
# the values previously held by the row before it in that cycle).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    8  = @{ D = 44617; L = "Primera"; M = 200; N = 6000;  O = 7000;  P = 6500;  S = 3250 }
    9  = @{ D = 44195; L = "Primera"; M = 200; N = 3000;  O = 3500;  P = 3250;  S = 1625 }
    10 = @{ D = 44195; L = "Segunda"; M = 100; N = 2500;  O = 2500;  P = 2500;  S = 1250 }
    11 = @{ D = 44532; L = "Primera"; M = 100; N = 10000; O = 10000; P = 10000; S = 5000 }
    12 = @{ D = 44532; L = "Segunda"; M = 100; N = 8000;  O = 8000;  P = 8000;  S = 4000 }
    15 = @{ D = 44559; L = "Primera"; M = 200; N = 6000;  O = 7000;  P = 6500;  S = 3250 }
    16 = @{ D = 44559; L = "Segunda"; M = 100; N = 5000;  O = 5000;  P = 5000;  S = 2500 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("S$r").Value = $vals.S
}
